$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting the existing "Brasil" row (and anything below) down.
$ws.Rows("8:8").Insert()

# Fill the newly inserted row 8 with the Sergipe data.
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = "Taxa de estupro"
$ws.Range("C8").Value = "31/12/2021"
$ws.Range("D8").Value = 33.73995178052012
$ws.Range("E8").Value = "14º"
